$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------------
# 1) Make room: the five "quarter marker" rows currently sitting in rows
#    244-248 need to end up in rows 249-253 (same values/format), because
#    rows 242-248 are about to be filled in with new leave-record data.
#    Copy formats first (this reproduces the exact style indices instead of
#    minting new ones), then re-create the values and the EARNED( ) helper
#    column formula in the new location. Do this before anything below
#    touches rows 244-248, since those rows are used as the format source.
# ---------------------------------------------------------------------------
$ws.Range("A244:K248").Copy()
$ws.Range("A249:K253").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A249").Value2 = 45261
$ws.Range("A251").Value2 = 45292
$ws.Range("A252").Value2 = 45323

# A250 is a text "2024" label (quote-prefixed in the source file) - assign
# it with a leading apostrophe so it is stored as text, not as the number
# 2024, then restore the exact quote-prefixed label style from row 245
# (still untouched at this point in the script).
$ws.Range("A250").Value = "'2024"
$ws.Range("A245").Copy()
$ws.Range("A250").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("G249").Formula = $formula
$ws.Range("G250").Formula = $formula
$ws.Range("G251").Formula = $formula
$ws.Range("G252").Formula = $formula
$ws.Range("G253").Formula = $formula

# ---------------------------------------------------------------------------
# 2) Fix up formatting quirks that don't come along with the new entries
#    below:
#    - row 245 used to be the "2024" label row (style 47 on col A); now it
#      is a regular data row, so col A reverts to the plain style (40),
#      copied from row 244 while that row's format is still the plain one.
#    - row 248 used to be the trailing footer row (styles 49/15/41/42/12);
#      now it is a regular data row, so its whole style set reverts to the
#      plain row pattern, copied from row 247.
#    Both source rows (244, 247) are still untouched at this point.
# ---------------------------------------------------------------------------
$ws.Range("A244").Copy()
$ws.Range("A245").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A247:K247").Copy()
$ws.Range("A248:K248").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------------
# 3) Row 240: EARNED (col C) now carries 1.25 (the EARNED( ) helper column G
#    recalculates automatically off the table formula already in place).
# ---------------------------------------------------------------------------
$ws.Range("C240").Value2 = 1.25

# ---------------------------------------------------------------------------
# 4) New leave entries in rows 242-248 (their old PERIOD marker values are
#    cleared out first since they have all been relocated to rows 249-253
#    above).
# ---------------------------------------------------------------------------
$ws.Range("A242").ClearContents()
$ws.Range("B242").Value = "SL(2-0-0)"
$ws.Range("H242").Value2 = 2
$ws.Range("K242").Value = "9/28 , 10/2/2023"

$ws.Range("A243").Value2 = 45200
$ws.Range("B243").Value = "VL(3-0-0)"
$ws.Range("C243").Value2 = 1.25
$ws.Range("D243").Value2 = 3
$ws.Range("K243").Value = "10/19,20,24/2023"

$ws.Range("A244").ClearContents()
$ws.Range("B244").Value = "SL(3-0-0)"
$ws.Range("H244").Value2 = 3
$ws.Range("K244").Value = "10/11-13/2023"

$ws.Range("A245").ClearContents()
$ws.Range("B245").Value = "VL(5-0-0)"
$ws.Range("D245").Value2 = 5
$ws.Range("K245").Value = "10/25-27,31 - 11/3/2023"

$ws.Range("A246").Value2 = 45231
$ws.Range("B246").Value = "SL(1-0-0)"
$ws.Range("C246").Value2 = 1.25
$ws.Range("H246").Value2 = 1
$ws.Range("K246").Value2 = 45239

$ws.Range("A247").ClearContents()
$ws.Range("B247").Value = "SL(1-0-0)"
$ws.Range("H247").Value2 = 1
$ws.Range("K247").Value2 = 45247

$ws.Range("B248").Value = "SL(1-0-0)"
$ws.Range("H248").Value2 = 1
$ws.Range("K248").Value2 = 45258

# K242, K246, K247, K248 pick up the date number format (style 48) even
# though K242's own content ends up being a free-text date string rather
# than a real date serial.
$ws.Range("K240").Copy()
$ws.Range("K242").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K246").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K247").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K248").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# re-assert the values (PasteSpecial formats-only shouldn't disturb them,
# but make sure they are exactly what we intend after the format touch-ups)
$ws.Range("K242").Value = "9/28 , 10/2/2023"
$ws.Range("K246").Value2 = 45239
$ws.Range("K247").Value2 = 45247
$ws.Range("K248").Value2 = 45258

# ---------------------------------------------------------------------------
# 5) Grow Table1 to cover the five extra rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K253"))

# ---------------------------------------------------------------------------
# 6) Leave the selection roughly where the editor left it.
# ---------------------------------------------------------------------------
$ws.Range("B249").Select()
$excel.ActiveWindow.ScrollRow = 242
